# edit.ps1 - reproduces the OOXML diff via PowerPoint COM interop
#
# Summary of changes applied:
#   1. Date placeholder field text "11/23/2024" -> "11/24/2024" on the
#      slide master and all 11 slide layouts.
#   2. The "NGINX LOAD BALANCER" textbox on slide 1 is resized/repositioned
#      and its text shortened to "LOAD BALANCER".

function Get-ShapeByName {
    param($Shapes, [string]$NamePrefix)
    for ($i = 1; $i -le $Shapes.Count; $i++) {
        $sh = $Shapes.Item($i)
        if ($sh.Name -like "$NamePrefix*") {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Refresh the cached "datetimeFigureOut" date text everywhere it is
#    rendered: the slide master plus every slide layout.
# ---------------------------------------------------------------------
$newDate = "11/24/2024"

$master = $p.SlideMaster

$masterDateSh = Get-ShapeByName $master.Shapes "Date Placeholder"
if ($masterDateSh -ne $null) {
    $masterDateSh.TextFrame.TextRange.Text = $newDate
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    $layoutDateSh = Get-ShapeByName $layout.Shapes "Date Placeholder"
    if ($layoutDateSh -ne $null) {
        $layoutDateSh.TextFrame.TextRange.Text = $newDate
    }
}

# ---------------------------------------------------------------------
# 2. Slide 1: move/resize the "NGINX LOAD BALANCER" textbox and rename
#    it to "LOAD BALANCER".
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)
$lbShape = Get-ShapeByName $s.Shapes "TextBox 4"

# EMU -> point helpers (1 pt = 12700 EMU). The PowerPoint object model
# stores Left/Top/Width/Height as single-precision points, so nudge the
# requested values to land on the exact target EMU after round-trip.
$lbShape.Left = 256.9085388370079
$lbShape.Width = 226.9196471992126
$lbShape.TextFrame.TextRange.Text = "LOAD BALANCER"
